# Update cryptocurrency price/volume data per the latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''70.785.59'
$ws.Range('E2').Value = '  -3.34%  '
$ws.Range('D3').Value = '''3.849.54'
$ws.Range('E3').Value = '  -3.66%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''597.96'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').Value = '''168.75'
$ws.Range('E6').Value = '  +3.05%  '
$ws.Range('D7').Value = '''0.667'
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = '''0.743'
$ws.Range('E9').Value = '  -1.09%  '
$ws.Range('E10').Value = '  +3.72%  '
$ws.Range('D11').Value = '''53.12'
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('D13').Value = '''11.32'
$ws.Range('E13').Value = '  +2.91%  '
$ws.Range('D14').Value = '''4.455.10'
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').Value = '''21.17'
$ws.Range('E15').Value = '  +3.67%  '
$ws.Range('D16').Value = '''3.855.36'
$ws.Range('E16').Value = '  -3.40%  '
$ws.Range('D17').Value = '''13.87'
$ws.Range('E17').Value = '  -2.02%  '
$ws.Range('E18').Value = '  -5.45%  '
$ws.Range('E19').Value = '  -2.40%  '
$ws.Range('D20').Value = '''70.586.30'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').Value = '''438.07'
$ws.Range('E21').Value = '  -0.24%  '
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '''94.50'
$ws.Range('E23').Value = '  -2.03%  '
$ws.Range('E24').Value = '  -4.83%  '
$ws.Range('D25').Value = '''13.81'
$ws.Range('E25').Value = '  -3.66%  '
$ws.Range('D26').Value = '''11.59'
$ws.Range('E26').Value = '  +2.34%  '
$ws.Range('E27').Value = '  -8.36%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +0.73%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '''8.47'
$ws.Range('E30').Value = '  +8.38%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''34.98'
$ws.Range('E31').Value = '  -3.74%  '
$ws.Range('D32').Value = '''13.50'
$ws.Range('E32').Value = '  -2.19%  '
$ws.Range('D33').Value = '''48.24'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '''0.126'
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').Value = '''68.90'
$ws.Range('E35').Value = '  -2.94%  '
$ws.Range('D36').Value = '''0.0₃0979'
$ws.Range('E36').Value = '  +8.28%  '
$ws.Range('D37').Value = '''635.67'
$ws.Range('E37').Value = '  -5.37%  '
$ws.Range('D38').Value = '''0.432'
$ws.Range('E38').Value = '  -1.23%  '
$ws.Range('D39').Value = '''0.146'
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.13%  '
$ws.Range('E41').Value = '  +0.10%  '
$ws.Range('D42').Value = '''3.23'
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('D43').Value = '''2.89'
$ws.Range('E43').Value = '  +9.81%  '
$ws.Range('D44').Value = '''3.15'
$ws.Range('E44').Value = '  +17.92%  '
$ws.Range('E45').Value = '  -4.61%  '
$ws.Range('D46').Value = '''10.05'
$ws.Range('E46').Value = '  -5.78%  '
$ws.Range('E47').Value = '  -4.10%  '
$ws.Range('B48').Value = 'WEMIXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D48').Value = '''2.88'
$ws.Range('E48').Value = '  -13.54%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '''2.906.95'
$ws.Range('E49').Value = '  +0.41%  '
$ws.Range('E50').Value = '  -3.95%  '
$ws.Range('E51').Value = '  +2.00%  '
